$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 28 for the new test mail ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A28").Value = "Wanneer zijn jullie open?"
$ws.Range("B28").Value = "mailmind.test@zohomail.eu"
$ws.Range("C28").Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Range("D28").Value = "Openingstijden / Locatie"
$ws.Range("F28").Value = "2025-06-26 22:54:55"
$ws.Range("G28").Value = "Nee"
$ws.Range("H28").Value = "Nee"
$ws.Range("I28").Value = "Nee"

# Extend the existing conditional formatting ranges so they keep covering
# the whole column down to the newly added row 28.
$ws.Range("D2:D27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D28"))
$ws.Range("G2:G27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G28"))
$ws.Range("H2:H27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H28"))
$ws.Range("I2:I27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I28"))

# --- Dashboard sheet: bump the "Openingstijden / Locatie" count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 5
